$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last data row (23) into new row 24, preserving all formatting.
$ws.Rows.Item(23).Copy()
$ws.Rows.Item(24).Insert()

# Fill in the values for the new "Event" / "Event with image" test row.
$ws.Cells.Item(24, 1).Value = "Event"
$ws.Cells.Item(24, 2).Value = 145278
$ws.Cells.Item(24, 3).Value = "Event with image"
$ws.Cells.Item(24, 4).Value = "New collection"
$ws.Cells.Item(24, 6).Value = "eProcurement"
$ws.Cells.Item(24, 7).Value = "Yes"
$ws.Cells.Item(24, 8).Value = "Yes"
$ws.Cells.Item(24, 13).Value = "doe@example.com"

# Re-create the mailto hyperlink on the Collection Owner cell, then restore the
# original cell formatting/value that adding the hyperlink disturbs.
$ws.Hyperlinks.Add($ws.Cells.Item(24, 13), "mailto:doe@example.com") | Out-Null
$ws.Cells.Item(23, 13).Copy($ws.Cells.Item(24, 13))
$ws.Cells.Item(24, 13).Value = "doe@example.com"

# Grow the table/autofilter range to include the new row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:Q24"))

# Keep the hidden _FilterDatabase defined name in sync with the table range.
$fd = $wb.Names.Item(1)
$fd.RefersTo = "='1. Content items'!`$A`$1:`$Q`$24"
